# Auto-generated edit script applying the Omega_Profits diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 2694.0417
$ws.Range("I19").Value = 2558
$ws.Range("J19").Value = 2809.1538
$ws.Range("K19").Value = 2558
$ws.Range("L19").Value = 2809.1538
$ws.Range("M19").Value = -2383
$ws.Range("N19").Value = -3159.1538
$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()
$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()
$ws.Range("H103").Value = 1679.4
$ws.Range("I103").Value = 1100
$ws.Range("J103").Value = 1927.7142
$ws.Range("K103").Value = 3300
$ws.Range("L103").Value = 5783.142599999999
$ws.Range("M103").Value = -2714
$ws.Range("N103").Value = -6955.142599999999
$ws.Range("H107").Value = 1314.8
$ws.Range("I107").Value = 1064.4166
$ws.Range("K107").Value = 1064.4166
$ws.Range("M107").Value = 855.5834
$ws.Range("H111").Value = 9077.857
$ws.Range("I111").Value = 3596.1428
$ws.Range("J111").Value = 14559.571
$ws.Range("K111").Value = 10788.4284
$ws.Range("L111").Value = 43678.713
$ws.Range("M111").Value = -7721.428400000001
$ws.Range("N111").Value = -49812.713
$ws.Range("H113").Value = 2113.6924
$ws.Range("J113").Value = 1969.4286
$ws.Range("L113").Value = 1969.4286
$ws.Range("N113").Value = -8477.428599999999
$ws.Range("H125").Value = 1875
$ws.Range("I125").Value = 1500
$ws.Range("K125").Value = 13500
$ws.Range("M125").Value = -11040
$ws.Range("H132").Value = 2247.7659
$ws.Range("I132").Value = 2253.152
$ws.Range("K132").Value = 6759.456
$ws.Range("M132").Value = -4229.456
$ws.Range("H138").Value = 2633.803
$ws.Range("I138").Value = 1474.4615
$ws.Range("K138").Value = 4423.3845
$ws.Range("M138").Value = 716.6154999999999
$ws.Range("H141").Value = 2843.2856
$ws.Range("I141").Value = 2709.0212
$ws.Range("J141").Value = 5998.5
$ws.Range("K141").Value = 8127.0636
$ws.Range("L141").Value = 17995.5
$ws.Range("M141").Value = -2947.0636
$ws.Range("N141").Value = -28355.5
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2770.8406
$ws.Range("I32").Value = 1169.95
$ws.Range("K32").Value = 1169.95
$ws.Range("M32").Value = -882.95
$ws.Range("H45").Value = 2440.4736
$ws.Range("I45").Value = 2322
$ws.Range("K45").Value = 2322
$ws.Range("M45").Value = -1945
$ws.Range("H61").Value = 4845.5
$ws.Range("I61").Value = 4254.4287
$ws.Range("J61").Value = 8983
$ws.Range("K61").Value = 4254.4287
$ws.Range("L61").Value = 8983
$ws.Range("M61").Value = -4042.4287
$ws.Range("N61").Value = -9407
$ws.Range("H63").Value = 2544
$ws.Range("I63").Value = 2602.6667
$ws.Range("J63").Value = 2500
$ws.Range("K63").Value = 2602.6667
$ws.Range("L63").Value = 2500
$ws.Range("M63").Value = -1916.6667
$ws.Range("N63").Value = -3872
$ws.Range("H66").Value = 2544
$ws.Range("I66").Value = 2602.6667
$ws.Range("J66").Value = 2500
$ws.Range("K66").Value = 13013.3335
$ws.Range("L66").Value = 12500
$ws.Range("M66").Value = -9581.333500000001
$ws.Range("N66").Value = -19364
$ws.Range("H74").Value = 1502.7742
$ws.Range("J74").Value = 831
$ws.Range("L74").Value = 831
$ws.Range("N74").Value = -2579
$ws.Range("H77").Value = 1502.7742
$ws.Range("J77").Value = 831
$ws.Range("L77").Value = 4155
$ws.Range("N77").Value = -12891
$ws.Range("H104").Value = 25225
$ws.Range("J104").Value = 25225
$ws.Range("L104").Value = 25225
$ws.Range("N104").Value = -32213
$ws.Range("H122").Value = 3585.0435
$ws.Range("I122").Value = 3283.1
$ws.Range("K122").Value = 9849.299999999999
$ws.Range("M122").Value = -7399.299999999999
$ws.Range("H132").Value = 3511.25
$ws.Range("I132").Value = 3346.2222
$ws.Range("K132").Value = 10038.6666
$ws.Range("M132").Value = -7508.6666
$ws.Range("H136").Value = 4845.5
$ws.Range("I136").Value = 4254.4287
$ws.Range("J136").Value = 8983
$ws.Range("K136").Value = 12763.2861
$ws.Range("L136").Value = 26949
$ws.Range("M136").Value = -10213.2861
$ws.Range("N136").Value = -32049
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 11137388
$ws.Range("J7").Value = 8338666
$ws.Range("L7").Value = 8338666
$ws.Range("N7").Value = -8338892
$ws.Range("H86").Value = 31486134
$ws.Range("I86").Value = 94445780
$ws.Range("J86").Value = 6312.5
$ws.Range("K86").Value = 94445780
$ws.Range("L86").Value = 6312.5
$ws.Range("M86").Value = -94444657
$ws.Range("N86").Value = -8558.5
$ws.Range("H89").Value = 31486134
$ws.Range("I89").Value = 94445780
$ws.Range("J89").Value = 6312.5
$ws.Range("K89").Value = 472228900
$ws.Range("L89").Value = 31562.5
$ws.Range("M89").Value = -472223284
$ws.Range("N89").Value = -42794.5
$ws.Range("H99").Value = 2293.2104
$ws.Range("I99").Value = 1998.4375
$ws.Range("J99").Value = 3865.3333
$ws.Range("K99").Value = 1998.4375
$ws.Range("L99").Value = 3865.3333
$ws.Range("M99").Value = -500.4375
$ws.Range("N99").Value = -6861.3333
$ws.Range("H107").Value = 2370.7144
$ws.Range("I107").Value = 1972.909
$ws.Range("J107").Value = 3829.3333
$ws.Range("K107").Value = 1972.909
$ws.Range("L107").Value = 3829.3333
$ws.Range("M107").Value = -52.90900000000011
$ws.Range("N107").Value = -7669.3333
$ws.Range("H134").Value = 3372.3928
$ws.Range("I134").Value = 3349.1482
$ws.Range("K134").Value = 10047.4446
$ws.Range("M134").Value = -7512.444600000001
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2772.1392
$ws.Range("I31").Value = 3297.5405
$ws.Range("K31").Value = 3297.5405
$ws.Range("M31").Value = -3002.5405
$ws.Range("H34").Value = 2772.1392
$ws.Range("I34").Value = 3297.5405
$ws.Range("K34").Value = 3297.5405
$ws.Range("M34").Value = -3095.5405
$ws.Range("H51").Value = 28299
$ws.Range("I51").Value = 28299
$ws.Range("K51").Value = 28299
$ws.Range("M51").Value = -27563
$ws.Range("H59").Value = 65713.42999999999
$ws.Range("J59").Value = 79998.8
$ws.Range("L59").Value = 79998.8
$ws.Range("N59").Value = -82288.8
$ws.Range("H61").Value = 28299
$ws.Range("I61").Value = 28299
$ws.Range("K61").Value = 28299
$ws.Range("M61").Value = -27951
$ws.Range("H86").Value = 17829.7
$ws.Range("I86").Value = 26159.4
$ws.Range("K86").Value = 26159.4
$ws.Range("M86").Value = -25036.4
$ws.Range("H89").Value = 17829.7
$ws.Range("I89").Value = 26159.4
$ws.Range("K89").Value = 130797
$ws.Range("M89").Value = -125181
$ws.Range("H132").Value = 5974.7144
$ws.Range("J132").Value = 1100
$ws.Range("L132").Value = 3300
$ws.Range("N132").Value = -8360
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 4377.636
$ws.Range("I3").Value = 3615.5
$ws.Range("K3").Value = 10846.5
$ws.Range("M3").Value = -10734.5
$ws.Range("H107").Value = 539.3
$ws.Range("I107").Value = 324.125
$ws.Range("K107").Value = 972.375
$ws.Range("M107").Value = 947.625
$ws.Range("H108").Value = 5026.5
$ws.Range("I108").Value = 5026.5
$ws.Range("K108").Value = 15079.5
$ws.Range("M108").Value = -12199.5
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 8545.799999999999
$ws.Range("I80").Value = 8247.333000000001
$ws.Range("J80").Value = 8993.5
$ws.Range("K80").Value = 8247.333000000001
$ws.Range("L80").Value = 8993.5
$ws.Range("M80").Value = -7249.333000000001
$ws.Range("N80").Value = -10989.5
$ws.Range("H83").Value = 8545.799999999999
$ws.Range("I83").Value = 8247.333000000001
$ws.Range("J83").Value = 8993.5
$ws.Range("K83").Value = 41236.665
$ws.Range("L83").Value = 44967.5
$ws.Range("M83").Value = -36244.665
$ws.Range("N83").Value = -54951.5
$ws.Range("H102").Value = 1396.5
$ws.Range("I102").Value = 1441.4546
$ws.Range("K102").Value = 1441.4546
$ws.Range("M102").Value = 180.5454
$ws.Range("H126").Value = 4649.5
$ws.Range("I126").Value = 2499.6667
$ws.Range("K126").Value = 7499.000100000001
$ws.Range("M126").Value = -5029.000100000001
$ws.Range("H132").Value = 6964.081
$ws.Range("I132").Value = 7014.147
$ws.Range("J132").Value = 6396.6665
$ws.Range("K132").Value = 21042.441
$ws.Range("L132").Value = 19189.9995
$ws.Range("M132").Value = -18512.441
$ws.Range("N132").Value = -24249.9995
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1277.7778
$ws.Range("J22").Value = 1625
$ws.Range("L22").Value = 1625
$ws.Range("N22").Value = -2215
$ws.Range("H27").Value = 1277.7778
$ws.Range("J27").Value = 1625
$ws.Range("L27").Value = 1625
$ws.Range("N27").Value = -1839
$ws.Range("H46").Value = 6898.4165
$ws.Range("I46").Value = 5830.3335
$ws.Range("J46").Value = 7966.5
$ws.Range("K46").Value = 5830.3335
$ws.Range("L46").Value = 7966.5
$ws.Range("M46").Value = -5642.3335
$ws.Range("N46").Value = -8342.5
$ws.Range("H55").Value = 957.63635
$ws.Range("I55").Value = 855.4286
$ws.Range("J55").Value = 1136.5
$ws.Range("K55").Value = 855.4286
$ws.Range("L55").Value = 1136.5
$ws.Range("M55").Value = -682.4286
$ws.Range("N55").Value = -1482.5
$ws.Range("H61").Value = 2745.158
$ws.Range("I61").Value = 2712.6924
$ws.Range("J61").Value = 2815.5
$ws.Range("K61").Value = 2712.6924
$ws.Range("L61").Value = 2815.5
$ws.Range("M61").Value = -2510.6924
$ws.Range("N61").Value = -3219.5
$ws.Range("H68").Value = 3327.875
$ws.Range("I68").Value = 2770.6667
$ws.Range("J68").Value = 4999.5
$ws.Range("K68").Value = 2770.6667
$ws.Range("L68").Value = 4999.5
$ws.Range("M68").Value = -2021.6667
$ws.Range("N68").Value = -6497.5
$ws.Range("H71").Value = 3327.875
$ws.Range("I71").Value = 2770.6667
$ws.Range("J71").Value = 4999.5
$ws.Range("K71").Value = 13853.3335
$ws.Range("L71").Value = 24997.5
$ws.Range("M71").Value = -10109.3335
$ws.Range("N71").Value = -32485.5
$ws.Range("H82").Value = 934.5
$ws.Range("I82").Value = 200
$ws.Range("J82").Value = 1179.3334
$ws.Range("K82").Value = 200
$ws.Range("L82").Value = 1179.3334
$ws.Range("M82").Value = 161
$ws.Range("N82").Value = -1901.3334
$ws.Range("H85").Value = 934.5
$ws.Range("I85").Value = 200
$ws.Range("J85").Value = 1179.3334
$ws.Range("K85").Value = 200
$ws.Range("L85").Value = 1179.3334
$ws.Range("M85").Value = 1048
$ws.Range("N85").Value = -3675.3334
$ws.Range("H113").Value = 2745.158
$ws.Range("I113").Value = 2712.6924
$ws.Range("J113").Value = 2815.5
$ws.Range("K113").Value = 2712.6924
$ws.Range("L113").Value = 2815.5
$ws.Range("M113").Value = -542.6923999999999
$ws.Range("N113").Value = -7155.5
$ws.Range("H122").Value = 9268.619000000001
$ws.Range("I122").Value = 9770.895
$ws.Range("K122").Value = 29312.685
$ws.Range("M122").Value = -26862.685
$ws.Range("H127").Value = 60499.75
$ws.Range("J127").Value = 60499.75
$ws.Range("L127").Value = 60499.75
$ws.Range("N127").Value = -70419.75
$ws.Range("H132").Value = 2864.04
$ws.Range("I132").Value = 2809.762
$ws.Range("J132").Value = 3149
$ws.Range("K132").Value = 8429.286
$ws.Range("L132").Value = 9447
$ws.Range("M132").Value = -5899.286
$ws.Range("N132").Value = -14507
$ws.Range("H139").Value = 89999
$ws.Range("J139").Value = 89999
$ws.Range("L139").Value = 89999
$ws.Range("N139").Value = -100279
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 6667262.5
$ws.Range("I2").Value = 894
$ws.Range("K2").Value = 894
$ws.Range("M2").Value = -782
$ws.Range("H39").Value = 20044
$ws.Range("I39").Value = 20044
$ws.Range("K39").Value = 20044
$ws.Range("M39").Value = -19631
$ws.Range("H45").Value = 47101
$ws.Range("J45").Value = 42686.75
$ws.Range("L45").Value = 42686.75
$ws.Range("N45").Value = -43668.75
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()
$ws.Range("H95").Value = 81343
$ws.Range("J95").Value = 81343
$ws.Range("L95").Value = 81343
$ws.Range("N95").Value = -86835
$ws.Range("H96").Value = 845.75
$ws.Range("I96").Value = 814.8
$ws.Range("J96").Value = 1000.5
$ws.Range("K96").Value = 814.8
$ws.Range("L96").Value = 1000.5
$ws.Range("M96").Value = 558.2
$ws.Range("N96").Value = -3746.5
$ws.Range("H132").Value = 3943.122
$ws.Range("I132").Value = 3082.0908
$ws.Range("K132").Value = 9246.2724
$ws.Range("M132").Value = -6716.2724
$ws.Range("H137").Value = 80611.625
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 80611.625
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 80611.625
$ws.Range("M137").ClearContents()
$ws.Range("N137").Value = -90811.625
